$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Expand the instrumentation amplifier bypass capacitor designators and quantity (row 4)
$ws.Range("A4").Value = "C5-C10"
$ws.Range("B4").Value = 6

# Update the 0.47uF capacitor voltage rating from 6.3V to 10V (row 3, C3-C4 bypass caps)
$ws.Range("D3").Value = "0402, Ceramic, X5R, 0.47uF, 10%, 10V"

# Update the Digikey short link for that part (row 3)
$ws.Range("G3").Value = "http://www.digikey.com/short/3tbjbm"

# Update the Digikey cart link (row 1, column J)
$ws.Range("J1").Value = "http://www.digikey.com/short/3tbj5t"

# Update the view state: select D4 and reset the scroll position
$ws.Range("D4").Select()
